$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain text (to match source formatting)
$textCells = @("D5", "D6", "D10", "D11", "D17", "D19", "D20", "D22", "D30", "D31", "D36", "D38", "D39", "D40", "D43", "D46", "D47", "D48", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '61.951.92'
$ws.Range('E2').Value = '  +4.61%  '
$ws.Range('D3').Value = '3.414.33'
$ws.Range('E3').Value = '  +3.29%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '578.11'
$ws.Range('E5').Value = '  +3.26%  '
$ws.Range('D6').Value = '138.62'
$ws.Range('E6').Value = '  +9.05%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = '3.415.27'
$ws.Range('E8').Value = '  +3.31%  '
$ws.Range('E9').Value = '  +2.19%  '
$ws.Range('D10').Value = '7.49'
$ws.Range('E10').Value = '  +2.12%  '
$ws.Range('D11').Value = '0.127'
$ws.Range('E11').Value = '  +10.06%  '
$ws.Range('E12').Value = '  +6.79%  '
$ws.Range('D13').Value = '3.999.33'
$ws.Range('E13').Value = '  +3.04%  '
$ws.Range('E14').Value = '  +1.85%  '
$ws.Range('E15').Value = '  +8.83%  '
$ws.Range('D16').Value = '3.415.10'
$ws.Range('E16').Value = '  +2.94%  '
$ws.Range('D17').Value = '25.49'
$ws.Range('E17').Value = '  +6.60%  '
$ws.Range('D18').Value = '61.967.51'
$ws.Range('E18').Value = '  +4.09%  '
$ws.Range('D19').Value = '14.18'
$ws.Range('E19').Value = '  +7.45%  '
$ws.Range('D20').Value = '5.90'
$ws.Range('E20').Value = '  +5.15%  '
$ws.Range('E21').Value = '  +8.00%  '
$ws.Range('D22').Value = '390.09'
$ws.Range('E22').Value = '  +11.79%  '
$ws.Range('E23').Value = '  +3.98%  '
$ws.Range('D24').Value = '3.552.14'
$ws.Range('E24').Value = '  +3.26%  '
$ws.Range('E25').Value = '  +18.93%  '
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('E27').Value = '  +4.67%  '
$ws.Range('E28').Value = '  +11.51%  '
$ws.Range('E29').Value = '  +5.71%  '
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.14%  '
$ws.Range('D31').Value = '8.30'
$ws.Range('E31').Value = '  +7.14%  '
$ws.Range('E32').Value = '  +7.20%  '
$ws.Range('E33').Value = '  +4.62%  '
$ws.Range('D34').Value = '3.446.27'
$ws.Range('E34').Value = '  +3.28%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').Value = '23.58'
$ws.Range('E36').Value = '  +4.05%  '
$ws.Range('E37').Value = '  +5.71%  '
$ws.Range('D38').Value = '7.01'
$ws.Range('E38').Value = '  +4.05%  '
$ws.Range('D39').Value = '1.57'
$ws.Range('E39').Value = '  +7.18%  '
$ws.Range('D40').Value = '162.82'
$ws.Range('E40').Value = '  +3.48%  '
$ws.Range('E41').Value = '  +6.75%  '
$ws.Range('E42').Value = '  +16.30%  '
$ws.Range('D43').Value = '0.791'
$ws.Range('E43').Value = '  +7.01%  '
$ws.Range('E44').Value = '  +7.10%  '
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').Value = '4.48'
$ws.Range('E46').Value = '  +4.88%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '25.12'
$ws.Range('E47').Value = '  +10.76%  '
$ws.Range('D48').Value = '41.75'
$ws.Range('E48').Value = '  +3.60%  '
$ws.Range('E49').Value = '  +4.40%  '
$ws.Range('D50').Value = '23.12'
$ws.Range('E50').Value = '  +7.00%  '
$ws.Range('D51').Value = '2.374.68'
$ws.Range('E51').Value = '  +9.96%  '
